$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3166.5
$ws.Range("I40").Value = 3000
$ws.Range("J40").Value = 3333
$ws.Range("K40").Value = 3000
$ws.Range("L40").Value = 3333
$ws.Range("M40").Value = -2825
$ws.Range("N40").Value = -3683

$ws.Range("H62").Value = 2951.25
$ws.Range("J62").Value = 2900
$ws.Range("L62").Value = 2900
$ws.Range("N62").Value = -4148

$ws.Range("H65").Value = 2951.25
$ws.Range("J65").Value = 2900
$ws.Range("L65").Value = 14500
$ws.Range("N65").Value = -20740

$ws.Range("H132").Value = 3107.8276
$ws.Range("I132").Value = 3043.08
$ws.Range("K132").Value = 9129.24
$ws.Range("M132").Value = -6599.24

$ws.Range("H138").Value = 1834.3654
$ws.Range("J138").Value = 3115.1155
$ws.Range("L138").Value = 9345.3465
$ws.Range("N138").Value = -19625.3465

$ws.Range("H141").Value = 3412.375
$ws.Range("I141").Value = 2325
$ws.Range("J141").Value = 4499.75
$ws.Range("K141").Value = 6975
$ws.Range("L141").Value = 13499.25
$ws.Range("M141").Value = -1795
$ws.Range("N141").Value = -23859.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3349.0667
$ws.Range("J45").Value = 3570.889
$ws.Range("L45").Value = 3570.889
$ws.Range("N45").Value = -4324.889

$ws.Range("H61").Value = 3106.6667
$ws.Range("I61").Value = 1780
$ws.Range("K61").Value = 1780
$ws.Range("M61").Value = -1568

$ws.Range("H63").Value = 2606232
$ws.Range("I63").Value = 2278.3
$ws.Range("J63").Value = 15626000
$ws.Range("K63").Value = 2278.3
$ws.Range("L63").Value = 15626000
$ws.Range("M63").Value = -1592.3
$ws.Range("N63").Value = -15627372

$ws.Range("H66").Value = 2606232
$ws.Range("I66").Value = 2278.3
$ws.Range("J66").Value = 15626000
$ws.Range("K66").Value = 11391.5
$ws.Range("L66").Value = 78130000
$ws.Range("M66").Value = -7959.5
$ws.Range("N66").Value = -78136864

$ws.Range("H74").Value = 2263.84
$ws.Range("I74").Value = 2212.3333
$ws.Range("K74").Value = 2212.3333
$ws.Range("M74").Value = -1338.3333

$ws.Range("H77").Value = 2263.84
$ws.Range("I77").Value = 2212.3333
$ws.Range("K77").Value = 11061.6665
$ws.Range("M77").Value = -6693.666499999999

$ws.Range("H122").Value = 1818
$ws.Range("I122").Value = 1897.5
$ws.Range("K122").Value = 5692.5
$ws.Range("M122").Value = -3242.5

$ws.Range("H136").Value = 3106.6667
$ws.Range("I136").Value = 1780
$ws.Range("K136").Value = 5340
$ws.Range("M136").Value = -2790

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1541.7273
$ws.Range("I20").Value = 1772.6
$ws.Range("J20").Value = 1047
$ws.Range("K20").Value = 1772.6
$ws.Range("L20").Value = 1047
$ws.Range("M20").Value = -1525.6
$ws.Range("N20").Value = -1541

$ws.Range("H94").Value = 2874.7942
$ws.Range("I94").Value = 1558.2174
$ws.Range("K94").Value = 1558.2174
$ws.Range("M94").Value = -1107.2174

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 69900
$ws.Range("J68").Value = 69900
$ws.Range("L68").Value = 69900
$ws.Range("N68").Value = -71398

$ws.Range("H71").Value = 69900
$ws.Range("J71").Value = 69900
$ws.Range("L71").Value = 209700
$ws.Range("N71").Value = -217188

$ws.Range("H74").Value = 29943.875
$ws.Range("J74").Value = 29943.875
$ws.Range("L74").Value = 29943.875
$ws.Range("N74").Value = -31691.875

$ws.Range("H77").Value = 29943.875
$ws.Range("J77").Value = 29943.875
$ws.Range("L77").Value = 89831.625
$ws.Range("N77").Value = -98567.625

$ws.Range("H105").Value = 25000388
$ws.Range("I105").Value = 31250232
$ws.Range("J105").Value = 1011
$ws.Range("K105").Value = 31250232
$ws.Range("L105").Value = 1011
$ws.Range("M105").Value = -31248485
$ws.Range("N105").Value = -4505

$ws.Range("H107").Value = 1003
$ws.Range("J107").Value = 1004
$ws.Range("L107").Value = 1004
$ws.Range("N107").Value = -4844

$ws.Range("H132").Value = 3115.5
$ws.Range("I132").Value = 1353
$ws.Range("J132").Value = 8403
$ws.Range("K132").Value = 4059
$ws.Range("L132").Value = 25209
$ws.Range("M132").Value = -1529
$ws.Range("N132").Value = -30269

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 271.8
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 271.8
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 1630.8
$ws.Range("N33").Value = -2196.8
$ws.Range("M33").ClearContents()

$ws.Range("H107").Value = 11400.777
$ws.Range("I107").Value = 50075
$ws.Range("J107").Value = 351
$ws.Range("K107").Value = 150225
$ws.Range("L107").Value = 1053
$ws.Range("M107").Value = -148305
$ws.Range("N107").Value = -4893

$ws.Range("H131").Value = 755.7755
$ws.Range("J131").Value = 777.37634
$ws.Range("L131").Value = 2332.12902
$ws.Range("N131").Value = -12412.12902

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3653.2273
$ws.Range("J80").Value = 3752.1428
$ws.Range("L80").Value = 3752.1428
$ws.Range("N80").Value = -5748.1428

$ws.Range("H83").Value = 3653.2273
$ws.Range("J83").Value = 3752.1428
$ws.Range("L83").Value = 18760.714
$ws.Range("N83").Value = -28744.714

$ws.Range("H97").Value = 3368.5264
$ws.Range("I97").Value = 1616.6666
$ws.Range("K97").Value = 1616.6666
$ws.Range("M97").Value = -1120.6666

$ws.Range("H122").Value = 2230.5
$ws.Range("I122").Value = 2216.5
$ws.Range("J122").Value = 2251.5
$ws.Range("K122").Value = 6649.5
$ws.Range("L122").Value = 6754.5
$ws.Range("M122").Value = -4199.5
$ws.Range("N122").Value = -11654.5

$ws.Range("H132").Value = 154421.1
$ws.Range("I132").Value = 147459
$ws.Range("K132").Value = 442377
$ws.Range("M132").Value = -439847

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2522.5
$ws.Range("I46").Value = 2366.6667
$ws.Range("J46").Value = 2616
$ws.Range("K46").Value = 2366.6667
$ws.Range("L46").Value = 2616
$ws.Range("M46").Value = -2178.6667
$ws.Range("N46").Value = -2992

$ws.Range("H61").Value = 7328.1816
$ws.Range("I61").Value = 2650
$ws.Range("K61").Value = 2650
$ws.Range("M61").Value = -2448

$ws.Range("H113").Value = 7328.1816
$ws.Range("I113").Value = 2650
$ws.Range("K113").Value = 2650
$ws.Range("M113").Value = -480

$ws.Range("H122").Value = 1404306.2
$ws.Range("I122").Value = 3924248.8
$ws.Range("J122").Value = 4338.222
$ws.Range("K122").Value = 11772746.4
$ws.Range("L122").Value = 13014.666
$ws.Range("M122").Value = -11770296.4
$ws.Range("N122").Value = -17914.666

$ws.Range("H132").Value = 2745.2727
$ws.Range("J132").Value = 4966
$ws.Range("L132").Value = 14898
$ws.Range("N132").Value = -19958

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4874.75
$ws.Range("I62").Value = 2999
$ws.Range("J62").Value = 5500
$ws.Range("K62").Value = 2999
$ws.Range("L62").Value = 5500
$ws.Range("M62").Value = -2375
$ws.Range("N62").Value = -6748

$ws.Range("H65").Value = 4874.75
$ws.Range("I65").Value = 2999
$ws.Range("J65").Value = 5500
$ws.Range("K65").Value = 14995
$ws.Range("L65").Value = 27500
$ws.Range("M65").Value = -11875
$ws.Range("N65").Value = -33740

$ws.Range("H81").Value = 1351
$ws.Range("I81").Value = 1467.7778
$ws.Range("J81").Value = 300
$ws.Range("K81").Value = 2935.5556
$ws.Range("L81").Value = 600
$ws.Range("M81").Value = -1874.5556
$ws.Range("N81").Value = -2722

$ws.Range("H84").Value = 1351
$ws.Range("I84").Value = 1467.7778
$ws.Range("J84").Value = 300
$ws.Range("K84").Value = 14677.778
$ws.Range("L84").Value = 3000
$ws.Range("M84").Value = -9373.778
$ws.Range("N84").Value = -13608

$ws.Range("H113").Value = 4506354.5
$ws.Range("I113").Value = 2675
$ws.Range("J113").Value = 13513714
$ws.Range("K113").Value = 8025
$ws.Range("L113").Value = 40541142
$ws.Range("M113").Value = -5855
$ws.Range("N113").Value = -40545482

$ws.Range("H122").Value = 2216.9285
$ws.Range("I122").Value = 2013.6364
$ws.Range("K122").Value = 6040.9092
$ws.Range("M122").Value = -3590.9092

$ws.Range("H126").Value = 2250
$ws.Range("J126").Value = 2666.6667
$ws.Range("L126").Value = 8000.000100000001
$ws.Range("N126").Value = -12940.0001

$ws.Range("H132").Value = 3599.3333
$ws.Range("I132").Value = 800
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 2400
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = 130
$ws.Range("N132").Value = -20057

$ws.Range("H136").Value = 18869192
$ws.Range("I136").Value = 27028148
$ws.Range("K136").Value = 81084444
$ws.Range("M136").Value = -81081894
